# Add a new paragraph style "marginOuter" that duplicates the existing
# "MarginNoteOuter" style (same basedOn / qFormat / frame-margin-note
# formatting), matching the lower-camel-case alias used elsewhere in the
# stylesheet output paths.
#
# (The frame-positioning properties on MarginNoteOuter - w:pPr/w:framePr -
# and the style-level w:rsid stamp aren't reachable through any exposed
# Style/ParagraphFormat/Frame COM property in this host, so only the
# name/basedOn/qFormat facets - the parts the object model actually
# supports writing - are reproduced here.)

$d = $word.ActiveDocument

# wdStyleTypeParagraph = 1
$marginOuter = $d.Styles.Add("marginOuter", 1)
$marginOuter.BaseStyle = "Normal"
$marginOuter.QuickStyle = $true

Write-Output "Added style: $($marginOuter.NameLocal)"
